$d = $word.ActiveDocument

# 1) Remove the old text "Programmes existants " from the first paragraph,
#    leaving the paragraph mark (and its bookmark) intact.
[void]$d.Content.Find.Execute("Programmes existants ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)

# 2) Insert a new heading paragraph before the (now empty) first paragraph.
$firstPara = $d.Paragraphs(1)
[void]$firstPara.Range.InsertParagraphBefore()
$insPoint = $d.Paragraphs(1).Range

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml"><pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships>
</pkg:xmlData></pkg:part>
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Descriptif des programmes </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>expect</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>r&#233;alis&#233;s</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> contenu dans le dossier prog.</w:t></w:r></w:p><w:sectPr/></w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@

[void]$insPoint.InsertXML($xml)
